$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - sheet1
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 13801
$wsExhibition.Range("F3").Value = 326
$wsExhibition.Range("F4").Value = 668
$wsExhibition.Range("F5").Value = 234
$wsExhibition.Range("F6").Value = 498
$wsExhibition.Range("F7").Value = 1420
$wsExhibition.Range("F8").Value = 135

# Sheet "全部类型" (All types) - sheet4
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 13801
$wsAll.Range("F3").Value = 326
$wsAll.Range("F4").Value = 668
$wsAll.Range("F5").Value = 234
$wsAll.Range("F8").Value = 498
$wsAll.Range("F9").Value = 1420
$wsAll.Range("F11").Value = 135
